$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text in C17 (shared string) to append the new content
$ws.Range("C17").Value = "Käyttäjän lisääminen ja login viimeistelty, Strava tokenien tallennus tietokantaan, yksinkertainen Strava aktiviteettien haku fronttiin"

# Update the hours value in B17 from 5 to 10
$ws.Range("B17").Value = 10

# Update row 17 height to 52
$ws.Rows.Item(17).RowHeight = 52

# Update the active cell/selection to B17
$ws.Range("B17").Select()
